$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was "M", now "B")
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9074889867841409
$ws.Range("C2").Value = 0.9581395348837209
$ws.Range("D2").Value = 0.9321266968325792
$ws.Range("E2").Value = 215

# Row 3 (was "B", now "M")
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.9217391304347826
$ws.Range("C3").Value = 0.8346456692913385
$ws.Range("D3").Value = 0.8760330578512396
$ws.Range("E3").Value = 127

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9146140586094618
$ws.Range("C5").Value = 0.8963926020875297
$ws.Range("D5").Value = 0.9040798773419094

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9127807067947592
$ws.Range("D6").Value = 0.9112966028248888
